# Fix the ASSISTS column (F) so that the values are stored as real numbers
# instead of text, and correct a handful of CHAMPION (H) names to "Milio".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F (ASSISTS): convert text "0"/"4"/"5"/... into numeric values ---
$assists = @{
    2  = 0;  3  = 0;  4  = 0;  5  = 0;  6  = 0;  7  = 0;  8  = 0;  9  = 0;
    10 = 0;  11 = 0;  12 = 0;  13 = 0;  14 = 0;  15 = 0;  16 = 0;  17 = 0;
    18 = 0;  19 = 0;  20 = 0;  21 = 4;  22 = 4;  23 = 4;  24 = 4;  25 = 5;
    26 = 5;  27 = 5;  28 = 6;  29 = 6;  30 = 6;  31 = 6;  32 = 6;  33 = 6;
    34 = 6;  35 = 6;  36 = 7;  37 = 8;  38 = 8;  39 = 8;  40 = 8;  41 = 8
}

foreach ($row in $assists.Keys) {
    $ws.Cells.Item($row, 6).Value = $assists[$row]
}

# --- Column H (CHAMPION): rename a few entries to "Milio" ---
$championRows = @(5, 11, 17, 23, 29, 35, 41)
foreach ($row in $championRows) {
    $ws.Cells.Item($row, 8).Value = "Milio"
}
